$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new rows (29, 30) to the bottom of the UNFI Bakery order sheet.
# All values are stored as text (matching the sheet's existing convention of
# t="inlineStr" cells for every column, including the numeric-looking
# Quantity/Cost columns) by prefixing numeric-looking literals with a
# leading apostrophe so Excel treats them as text instead of numbers, then
# resetting the cell style back to Normal so no stray formatting survives.

$ws.Range("A29").Value = "183096-7"
$ws.Range("B29").Value = "Clio - Greek Yogurt Bar Strawberry"
$ws.Range("C29").Value = "'1"
$ws.Range("D29").Value = "'15.45"
$ws.Range("E29").Value = "'15.45"

$ws.Range("A30").Value = "183090-0"
$ws.Range("B30").Value = "Clio - Greek Yogurt Bar Vanilla"
$ws.Range("C30").Value = "'1"
$ws.Range("D30").Value = "'15.45"
$ws.Range("E30").Value = "'15.45"

$ws.Range("A29:E30").Style = "Normal"
